# Applies the "corrections up into the probabilistic aspects" edit:
#  1. Refresh the cached date/time field text (datetimeFigureOut) from
#     3/16/2023 (en-US) / 16/03/2023 (en-DE) to 6/26/2023 / 26/06/2023
#     across the Notes Master, Slide Master, and every Slide Layout.
#  2. Correct the slide 1 textbox wording "FN or delayed TP?" ->
#     "FN and/or delayed TP?".

$p = $ppt.ActivePresentation

# --- 1. Notes Master: Date Placeholder (en-US, M/D/YYYY style) ---
$nm = $p.NotesMaster
for ($i = 1; $i -le $nm.Shapes.Count; $i++) {
    $sh = $nm.Shapes.Item($i)
    if ($sh.Name -like "Date Placeholder*") {
        $sh.TextFrame.TextRange.Text = "6/26/2023"
    }
}

# --- 2. Slide Master: Date Placeholder (en-DE, D/M/YYYY style) ---
$sm = $p.SlideMaster
for ($i = 1; $i -le $sm.Shapes.Count; $i++) {
    $sh = $sm.Shapes.Item($i)
    if ($sh.Name -like "Date Placeholder*") {
        $sh.TextFrame.TextRange.Text = "26/06/2023"
    }
}

# --- 3. Every Slide Layout: Date Placeholder (en-DE, D/M/YYYY style) ---
for ($j = 1; $j -le $sm.CustomLayouts.Count; $j++) {
    $cl = $sm.CustomLayouts.Item($j)
    for ($i = 1; $i -le $cl.Shapes.Count; $i++) {
        $sh = $cl.Shapes.Item($i)
        if ($sh.Name -like "Date Placeholder*") {
            $sh.TextFrame.TextRange.Text = "26/06/2023"
        }
    }
}

# --- 4. Slide 1 textbox wording correction ---
$s1 = $p.Slides.Item(1)
for ($i = 1; $i -le $s1.Shapes.Count; $i++) {
    $sh = $s1.Shapes.Item($i)
    if ($sh.HasTextFrame -and $sh.TextFrame.HasText) {
        if ($sh.TextFrame.TextRange.Text -eq "FN or delayed TP?") {
            $sh.TextFrame.TextRange.Text = "FN and/or delayed TP?"
        }
    }
}
